$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 448.5
$ws.Range("I28").Value = 466.6875
$ws.Range("K28").Value = 466.6875
$ws.Range("M28").Value = 18.3125

$ws.Range("H51").Value = 4250
$ws.Range("J51").Value = 4500
$ws.Range("L51").Value = 4500
$ws.Range("N51").Value = -5468

$ws.Range("H53").Value = 234.27272
$ws.Range("I53").Value = 152
$ws.Range("J53").Value = 302.83334
$ws.Range("K53").Value = 152
$ws.Range("L53").Value = 302.83334
$ws.Range("M53").Value = 485
$ws.Range("N53").Value = -1576.83334

$ws.Range("H76").Value = 4582.5
$ws.Range("I76").Value = 3875
$ws.Range("J76").Value = 5997.5
$ws.Range("K76").Value = 3875
$ws.Range("L76").Value = 5997.5
$ws.Range("M76").Value = -3560
$ws.Range("N76").Value = -6627.5

$ws.Range("H79").Value = 4582.5
$ws.Range("I79").Value = 3875
$ws.Range("J79").Value = 5997.5
$ws.Range("K79").Value = 3875
$ws.Range("L79").Value = 5997.5
$ws.Range("M79").Value = -2783
$ws.Range("N79").Value = -8181.5

$ws.Range("H86").Value = 8245.637000000001
$ws.Range("J86").Value = 6101
$ws.Range("L86").Value = 6101
$ws.Range("N86").Value = -8347

$ws.Range("H89").Value = 8245.637000000001
$ws.Range("J89").Value = 6101
$ws.Range("L89").Value = 30505
$ws.Range("N89").Value = -41737

$ws.Range("H92").Value = 1488.4
$ws.Range("I92").Value = 1497.5
$ws.Range("J92").Value = 1482.3334
$ws.Range("K92").Value = 1497.5
$ws.Range("L92").Value = 1482.3334
$ws.Range("M92").Value = -249.5
$ws.Range("N92").Value = -3978.3334

$ws.Range("H98").Value = 845.5
$ws.Range("I98").Value = 844.9677
$ws.Range("K98").Value = 844.9677
$ws.Range("M98").Value = 653.0323

$ws.Range("H106").Value = 3334.1667
$ws.Range("I106").Value = 3132.5
$ws.Range("K106").Value = 3132.5
$ws.Range("M106").Value = -2501.5

$ws.Range("H107").Value = 5210.1577
$ws.Range("I107").Value = 4362.625
$ws.Range("K107").Value = 4362.625
$ws.Range("M107").Value = -2442.625

$ws.Range("H122").Value = 845.5
$ws.Range("I122").Value = 844.9677
$ws.Range("K122").Value = 2534.9031
$ws.Range("M122").Value = -84.90309999999999

$ws.Range("H138").Value = 2155.2466
$ws.Range("I138").Value = 1088.4193
$ws.Range("K138").Value = 3265.2579
$ws.Range("M138").Value = 1874.7421

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6339.143
$ws.Range("I32").Value = 5749.855
$ws.Range("K32").Value = 5749.855
$ws.Range("M32").Value = -5462.855

$ws.Range("H43").Value = 18000
$ws.Range("I43").Value = 18000
$ws.Range("K43").Value = 18000
$ws.Range("M43").Value = -17687

$ws.Range("H97").Value = 128708.125
$ws.Range("I97").Value = 4077.5
$ws.Range("J97").Value = 502600
$ws.Range("K97").Value = 4077.5
$ws.Range("L97").Value = 502600
$ws.Range("M97").Value = -3581.5
$ws.Range("N97").Value = -503592

$ws.Range("H122").Value = 40838.777
$ws.Range("I122").Value = 5138
$ws.Range("K122").Value = 15414
$ws.Range("M122").Value = -12964

$ws.Range("H124").Value = 39664
$ws.Range("J124").Value = 39664
$ws.Range("L124").Value = 39664
$ws.Range("N124").Value = -49484

$ws.Range("H130").Value = 72998.75
$ws.Range("J130").Value = 72998.75
$ws.Range("L130").Value = 72998.75
$ws.Range("N130").Value = -83038.75

$ws.Range("H132").Value = 2391.8
$ws.Range("I132").Value = 2260.1853
$ws.Range("J132").Value = 3576.3333
$ws.Range("K132").Value = 6780.5559
$ws.Range("L132").Value = 10728.9999
$ws.Range("M132").Value = -4250.5559
$ws.Range("N132").Value = -15788.9999

$ws.Range("H134").Value = 114844
$ws.Range("J134").Value = 114844
$ws.Range("L134").Value = 114844
$ws.Range("N134").Value = -124984

$ws.Range("H139").Value = 151819.5
$ws.Range("J139").Value = 151819.5
$ws.Range("L139").Value = 151819.5
$ws.Range("N139").Value = -162099.5

$ws.Range("H140").Value = 129999.664
$ws.Range("J140").Value = 129999.664
$ws.Range("L140").Value = 129999.664
$ws.Range("N140").Value = -140359.664

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1306.6666
$ws.Range("I94").Value = 900
$ws.Range("J94").Value = 1510
$ws.Range("K94").Value = 900
$ws.Range("L94").Value = 1510
$ws.Range("M94").Value = -449
$ws.Range("N94").Value = -2412

$ws.Range("H126").Value = 77890
$ws.Range("J126").Value = 77890
$ws.Range("L126").Value = 77890
$ws.Range("N126").Value = -87770

$ws.Range("H134").Value = 2571.7
$ws.Range("I134").Value = 2475.5557
$ws.Range("K134").Value = 7426.6671
$ws.Range("M134").Value = -4891.6671

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2076.7778
$ws.Range("I58").Value = 958.4
$ws.Range("K58").Value = 958.4
$ws.Range("M58").Value = -755.4

$ws.Range("H86").Value = 4765118
$ws.Range("I86").Value = 8336220
$ws.Range("J86").Value = 3649.5
$ws.Range("K86").Value = 8336220
$ws.Range("L86").Value = 3649.5
$ws.Range("M86").Value = -8335097
$ws.Range("N86").Value = -5895.5

$ws.Range("H89").Value = 4765118
$ws.Range("I89").Value = 8336220
$ws.Range("J89").Value = 3649.5
$ws.Range("K89").Value = 41681100
$ws.Range("L89").Value = 18247.5
$ws.Range("M89").Value = -41675484
$ws.Range("N89").Value = -29479.5

$ws.Range("H107").Value = 1126.4482
$ws.Range("I107").Value = 1010.7619
$ws.Range("K107").Value = 1010.7619
$ws.Range("M107").Value = 909.2381

$ws.Range("H134").Value = 2516.509
$ws.Range("I134").Value = 2144.5652
$ws.Range("J134").Value = 4417.5557
$ws.Range("K134").Value = 6433.6956
$ws.Range("L134").Value = 13252.6671
$ws.Range("M134").Value = -3898.6956
$ws.Range("N134").Value = -18322.6671

$ws.Range("H136").Value = 2076.7778
$ws.Range("I136").Value = 958.4
$ws.Range("K136").Value = 2875.2
$ws.Range("M136").Value = -325.1999999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 10901.632
$ws.Range("I56").Value = 10901.632
$ws.Range("K56").Value = 10901.632
$ws.Range("M56").Value = -10371.632

$ws.Range("H87").Value = 7499.5
$ws.Range("I87").Value = 7499.5
$ws.Range("K87").Value = 22498.5
$ws.Range("M87").Value = -21250.5

$ws.Range("H90").Value = 7499.5
$ws.Range("I90").Value = 7499.5
$ws.Range("K90").Value = 67495.5
$ws.Range("M90").Value = -61255.5

$ws.Range("H118").Value = 6800
$ws.Range("I118").Value = 6800
$ws.Range("K118").Value = 20400
$ws.Range("M118").Value = -19157

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H116").Value = 105185.5
$ws.Range("I116").Value = 60000
$ws.Range("J116").Value = 120247.336
$ws.Range("K116").Value = 60000
$ws.Range("L116").Value = 120247.336
$ws.Range("M116").Value = -55411
$ws.Range("N116").Value = -129425.336

$ws.Range("H122").Value = 718411.9399999999
$ws.Range("J122").Value = 4475.6
$ws.Range("L122").Value = 13426.8
$ws.Range("N122").Value = -18326.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 18801.834
$ws.Range("I61").Value = 18801.834
$ws.Range("K61").Value = 18801.834
$ws.Range("M61").Value = -18599.834

$ws.Range("H113").Value = 18801.834
$ws.Range("I113").Value = 18801.834
$ws.Range("K113").Value = 18801.834
$ws.Range("M113").Value = -16631.834

$ws.Range("H128").Value = 10000
$ws.Range("J128").Value = 10000
$ws.Range("L128").Value = 10000
$ws.Range("N128").Value = -19960

$ws.Range("H134").Value = 116620.5
$ws.Range("J134").Value = 116620.5
$ws.Range("L134").Value = 116620.5
$ws.Range("N134").Value = -126760.5

$ws.Range("H137").Value = 89998.5
$ws.Range("J137").Value = 109997
$ws.Range("L137").Value = 109997
$ws.Range("N137").Value = -120197

$ws.Range("H138").Value = 98646.664
$ws.Range("J138").Value = 98646.664
$ws.Range("L138").Value = 98646.664
$ws.Range("N138").Value = -108926.664

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 7525.7334
$ws.Range("I81").Value = 6677.4
$ws.Range("J81").Value = 7949.9
$ws.Range("K81").Value = 13354.8
$ws.Range("L81").Value = 15899.8
$ws.Range("M81").Value = -12293.8
$ws.Range("N81").Value = -18021.8

$ws.Range("H84").Value = 7525.7334
$ws.Range("I84").Value = 6677.4
$ws.Range("J84").Value = 7949.9
$ws.Range("K84").Value = 66774
$ws.Range("L84").Value = 79499
$ws.Range("M84").Value = -61470
$ws.Range("N84").Value = -90107

$ws.Range("H100").Value = 1060.8572
$ws.Range("I100").Value = 1071.1666
$ws.Range("J100").Value = 999
$ws.Range("K100").Value = 2142.3332
$ws.Range("L100").Value = 1998
$ws.Range("M100").Value = -1601.3332
$ws.Range("N100").Value = -3080

$ws.Range("H113").Value = 3812.5715
$ws.Range("I113").Value = 3439.6
$ws.Range("K113").Value = 10318.8
$ws.Range("M113").Value = -8148.799999999999

$ws.Range("H123").Value = 85375.25
$ws.Range("J123").Value = 85375.25
$ws.Range("L123").Value = 85375.25
$ws.Range("N123").Value = -95175.25

$ws.Range("H127").Value = 79000
$ws.Range("J127").Value = 79000
$ws.Range("L127").Value = 79000
$ws.Range("N127").Value = -88920

$ws.Range("H140").Value = 108243.25
$ws.Range("J140").Value = 108243.25
$ws.Range("L140").Value = 108243.25
$ws.Range("N140").Value = -118603.25
